# Adds two new slides (8) next steps / a blank follow-up slide) to the end
# of the deck, matching the "Title and Content" layout already used by the
# rest of the presentation.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 15: "8) next steps"
# ---------------------------------------------------------------------
$idx15 = $p.Slides.Count + 1
$s15 = $p.Slides.Add($idx15, 2)

$title15 = $s15.Shapes.Item(1).TextFrame.TextRange
$title15.Text = "8) next steps"
$title15.LanguageID = "pt-BR"

$body15 = $s15.Shapes.Item(2).TextFrame.TextRange
$body15.Text = "Improve monitoring stages in pipeline`r" + `
    "`r" + `
    "Correlate complaints and delivery events using parcel_id or purchase_order_id whenever possible`r" + `
    "`r" + `
    "Create a crosswalk table allowing a more precise conection between delivery events and complaints`r" + `
    "`r" + `
    "Draw insights from drivers, recipients and stations in terms of number of deliveries and number of complaints`r" + `
    "`r" + `
    "Use LLM to evaluate the claims contents`r"
$body15.LanguageID = "pt-BR"

# Italicize "parcel_id " (including the trailing space) and
# "purchase_order_id" inside the third paragraph.
$fullText = $body15.Text
$parcelStart = $fullText.IndexOf("parcel_id")
$parcelRange = $body15.Characters($parcelStart + 1, 10)
$parcelRange.Font.Italic = $true

$purchaseStart = $fullText.IndexOf("purchase_order_id")
$purchaseRange = $body15.Characters($purchaseStart + 1, 17)
$purchaseRange.Font.Italic = $true

# Shrink text to fit the placeholder (as happens automatically once the
# long bullet list overflows the content box in real PowerPoint).
$bodyFrame15 = $s15.Shapes.Item(2).TextFrame
$bodyFrame15.AutoSize = 2

# ---------------------------------------------------------------------
# Slide 16: blank "Title and Content" slide
# ---------------------------------------------------------------------
$idx16 = $p.Slides.Count + 1
$s16 = $p.Slides.Add($idx16, 2)

$title16 = $s16.Shapes.Item(1).TextFrame.TextRange
$title16.LanguageID = "pt-BR"

$body16 = $s16.Shapes.Item(2).TextFrame.TextRange
$body16.LanguageID = "pt-BR"
